$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# New synonym rows (A = old term, B = new/canonical term)
$newRows = @(
    @("doubly toothed",   "doubly-toothed"),
    @("double-toothed",   "doubly-toothed"),
    @("subacuminate",     "acuminate"),
    @("sub-acuminate",    "acuminate"),
    @("coarsely toothed", "dentate"),
    @("finely toothed",   "serrate")
)

$startRow = 56
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
}

# Widen column B to fit the new content
$ws.Columns.Item(2).ColumnWidth = 17.666666666666668

# Update the view: scroll/active cell moved as a result of the new rows
$ws.Activate()
$ws.Range("B62").Select()
